$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "Rango" label to K1 ---
$ws.Range("K1").Value = "Rango"

# --- Rename shared text used by column I header ("Varianza" -> "Desviación") ---
$ws.Range("I1").Value = "Desviación"

# --- Row 2 edits ---
$ws.Range("H2").Value = 4
$ws.Range("K2").Value = 1

# --- Row 3 edits ---
$ws.Range("C3").Value = 50
$ws.Range("H3").Value = 0.2
$ws.Range("I3").Value = 0.2

# --- Row 4: clear all data (was Jaguar "Depredador Diurno" entry) ---
$ws.Range("A4:K4").ClearContents()

# --- Row 5: clear the leading identifying columns (A:E) ---
$ws.Range("A5:E5").ClearContents()

# --- View state: update the active selection to match the saved workbook view ---
$ws.Activate()
$ws.Range("G6").Select()
